$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 44348
$ws.Range("J3").Value = 150
$ws.Range("K3").Value = 7000
$ws.Range("L3").Value = 7000
$ws.Range("M3").Value = 7000
$ws.Range("O3").Value = 'Región del Maule'
$ws.Range("P3").Value = 194
$ws.Range("D4").Value = 44364
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 7000
$ws.Range("L4").Value = 7000
$ws.Range("M4").Value = 7000
$ws.Range("P4").Value = 194
$ws.Range("D5").Value = 44376
$ws.Range("K5").Value = 6500
$ws.Range("L5").Value = 6500
$ws.Range("M5").Value = 6500
$ws.Range("P5").Value = 181
$ws.Range("D6").Value = 44354
$ws.Range("J6").Value = 150
$ws.Range("N6").Value = '$/caja 36 atados'
$ws.Range("O6").Value = 'Región del Maule'
$ws.Range("P6").Value = 194
$ws.Range("Q6").Value = 36
$ws.Range("D7").Value = 44369
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 7000
$ws.Range("L7").Value = 7000
$ws.Range("M7").Value = 7000
$ws.Range("N7").Value = '$/caja 20 docenas'
$ws.Range("P7").Value = 7000
$ws.Range("Q7").Value = 1
$ws.Range("D9").Value = 44362
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 6500
$ws.Range("L9").Value = 6500
$ws.Range("M9").Value = 6500
$ws.Range("O9").Value = 'Región Metropolitana'
$ws.Range("P9").Value = 181
$ws.Range("D10").Value = 44355
$ws.Range("O10").Value = 'Región Metropolitana'
$ws.Range("D11").Value = 44342
$ws.Range("D12").Value = 44372
$ws.Range("J12").Value = 150
$ws.Range("D13").Value = 44371
$ws.Range("K13").Value = 6500
$ws.Range("L13").Value = 6500
$ws.Range("M13").Value = 6500
$ws.Range("P13").Value = 181
$ws.Range("D14").Value = 44386
$ws.Range("J14").Value = 200
$ws.Range("D15").Value = 44358
